# Correct encoder ticks count for the full rotation
#
# Adds a new "Servo" worksheet (with motor rpm/rps numbers) after the
# existing sheets, and makes it the active/visible sheet.

$wb = $excel.ActiveWorkbook

# Add the new worksheet right after the last existing sheet
# (NFCv2 LED driver, NFCv2 USART1, LM43603, GenericExcercise, Servo)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Servo"

# Populate the cells. Write C2 ("rpm @50Hz") before A2 ("motor") so that
# new entries land in the shared-string table in the expected order:
# rpm @50Hz, motor, rps @50Hz.
$ws.Range("C2").Value = "rpm @50Hz"
$ws.Range("A2").Value = "motor"
$ws.Range("B2").Value = 1395
$ws.Range("B3").Formula = "=B2/60"
$ws.Range("C3").Value = "rps @50Hz"

# Match the selection left on the sheet and make it the active tab
$ws.Range("C4").Select()
$ws.Activate()
